$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 170, shifting existing rows 170-243 down to 171-244
$ws.Rows.Item(170).Insert()

# Populate the newly inserted row 170 with data
$ws.Cells.Item(170, 1).Value = 11
$ws.Cells.Item(170, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(170, 3).Value = "Bíobío"
$ws.Cells.Item(170, 4).Value = 45205
$ws.Cells.Item(170, 5).Value = 8
$ws.Cells.Item(170, 6).Value = 100112043
$ws.Cells.Item(170, 7).Value = "Pepino ensalada"
$ws.Cells.Item(170, 8).Value = "Sin especificar"
$ws.Cells.Item(170, 9).Value = "Primera"
$ws.Cells.Item(170, 10).Value = 100
$ws.Cells.Item(170, 11).Value = 16000
$ws.Cells.Item(170, 12).Value = 17000
$ws.Cells.Item(170, 13).Value = 16500
$ws.Cells.Item(170, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(170, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(170, 16).Value = 275
$ws.Cells.Item(170, 17).Value = 60
$ws.Cells.Item(170, 18).Value = "Hortaliza"
